# Phase 2 sheet addition, modeled after the commit "Added stuff for phase 2".
# Strategy: duplicate "Phase 1" (Copy keeps formulas/styles/conditional
# formatting identical), rename the duplicate "Phase 2", wire up the
# sheet-scoped defined name "prevWBS" for it (same pattern as sheet 1),
# then overwrite the task rows with the Phase-2 content and extend the
# table with a few new rows.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Phase 1")

# --- 1. Duplicate the sheet and rename it -------------------------------
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Phase 2"

# --- 2. Sheet-scoped defined name, same shape as 'Phase 1'!prevWBS ------
$ws2.Names.Add("prevWBS", "='Phase 2'!`$A1048576")

# --- 3. The WBS-numbering helper column (A) is blank on the new sheet --
$ws2.Range("A4:A9").ClearContents()

# --- 4. Re-point the task table at the Phase 2 content ------------------
$ws2.Range("B4").Value = "Phase 2: EDA and baseline pipeline. Team Lead: Paul Miller"
$ws2.Range("C4").Value = "Team member"

$ws2.Range("B5").Value = "Abstract, organize notebook"
$ws2.Range("C5").Value = "Paul Miller"

$ws2.Range("B6").Value = "Load data"
$ws2.Range("C6").Value = "Glen Colletti"

$ws2.Range("B7").Value = "EDA "
$ws2.Range("C7").Value = "Alex Bordanca"

$ws2.Range("B8").Value = "Visual EDA"
$ws2.Range("C8").Value = "Alex Bordanca"

$ws2.Range("B9").Value = "Baseline models and pipelines. XGBoost, KNN, Logistic Regression"
$ws2.Range("C9").Value = "Glen Colletti"

$ws2.Range("B10").Value = "Create presentation slides"
$ws2.Range("C10").Value = "Glen Colletti"

$ws2.Range("B11").Value = "Credit Assignment"
$ws2.Range("C11").Value = "Paul Miller"

$ws2.Range("B12").Value = "Record video"
$ws2.Range("C12").Value = "All members"

# --- 5. New rows 9-12 get a lighter Arial font + soft gray top/bottom ---
#        borders (distinct from the original Aptos Narrow/gray22 styling
#        of rows 5-8, matching the new look-and-feel added for Phase 2).
$newRows = $ws2.Range("B9:C12")
$newRows.Font.Name = "Arial"
$newRows.Font.Size = 9
$newRows.Borders.Item(8).LineStyle = 1
$newRows.Borders.Item(8).Color = 12632256
$newRows.Borders.Item(9).LineStyle = 1
$newRows.Borders.Item(9).Color = 12632256

$existingRows = $ws2.Range("B5:C8")
$existingRows.Font.Name = "Arial"
$existingRows.Font.Size = 9

# C12 sits below the original bordered block, so it keeps no border.
$ws2.Range("C12").Borders.Item(8).LineStyle = -4142
$ws2.Range("C12").Borders.Item(9).LineStyle = -4142

# --- 6. View bookkeeping: Phase 2 becomes the active/selected tab -------
$ws1.Range("B15").Select()
$ws2.Activate()
$ws2.Range("F13").Select()
